$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.083.28"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "3.445.99"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'576.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.39%  "

$ws.Range("D6").Value = "'161.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.71%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.449.41"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("D9").Value = "'0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.95%  "

$ws.Range("D11").Value = "'0.125"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").Value = "4.035.30"
$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").Value = "'0.0000193"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "'28.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").Value = "65.033.03"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").Value = "3.416.11"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "'6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "

$ws.Range("D20").Value = "'14.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").Value = "'386.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.90%  "

$ws.Range("D22").Value = "'8.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.87%  "

$ws.Range("D23").Value = "'73.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.07%  "

$ws.Range("D24").Value = "'0.543"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value = "'0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.42%  "

$ws.Range("D27").Value = "'9.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").Value = "'6.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.42%  "

$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("E32").Value = "  -0.83%  "

$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "'6.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'23.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'7.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.58%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'161.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.64%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'1.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.89%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.033.19"
$ws.Range("E39").Value = "  +4.65%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0766"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.06%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'27.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.05%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.36%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'42.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.43%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0317"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.770"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.30%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'24.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.28%  "

$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "

$ws.Range("D49").Value = "'0.870"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.67%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.59%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'306.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.08%  "
